$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.061.88"
$ws.Range("E2").Value = "'  +0.98%  "
$ws.Range("D3").Value = "'1.847.66"
$ws.Range("E3").Value = "'  +2.15%  "
$ws.Range("E4").Value = "'  +0.18%  "
$ws.Range("D5").Value = "'233.67"
$ws.Range("E5").Value = "'  +0.52%  "
$ws.Range("D6").Value = "'0.621"
$ws.Range("E6").Value = "'  +2.91%  "
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("D8").Value = "'41.50"
$ws.Range("E8").Value = "'  +5.68%  "
$ws.Range("E9").Value = "'  +2.19%  "
$ws.Range("E10").Value = "'  +2.01%  "
$ws.Range("D11").Value = "'0.0982"
$ws.Range("E11").Value = "'  -1.08%  "
$ws.Range("D12").Value = "'2.111.92"
$ws.Range("E12").Value = "'  +1.99%  "
$ws.Range("D13").Value = "'11.55"
$ws.Range("E13").Value = "'  +4.57%  "
$ws.Range("D14").Value = "'1.850.65"
$ws.Range("E14").Value = "'  +2.25%  "
$ws.Range("D15").Value = "'0.675"
$ws.Range("E15").Value = "'  +0.95%  "
$ws.Range("D16").Value = "'4.70"
$ws.Range("E16").Value = "'  +2.67%  "
$ws.Range("D17").Value = "'35.070.06"
$ws.Range("E17").Value = "'  +1.02%  "
$ws.Range("D18").Value = "'70.11"
$ws.Range("E18").Value = "'  +0.77%  "
$ws.Range("D19").Value = "'0.0₃0791"
$ws.Range("E19").Value = "'  +0.73%  "
$ws.Range("D20").Value = "'240.75"
$ws.Range("E20").Value = "'  +0.47%  "
$ws.Range("D21").Value = "'12.21"
$ws.Range("E21").Value = "'  +2.36%  "
$ws.Range("D22").Value = "'4.78"
$ws.Range("E22").Value = "'  +3.00%  "
$ws.Range("E23").Value = "'  +0.12%  "
$ws.Range("E24").Value = "'  +3.65%  "
$ws.Range("D25").Value = "'171.61"
$ws.Range("E25").Value = "'  -0.24%  "
$ws.Range("D26").Value = "'7.90"
$ws.Range("E26").Value = "'  +2.30%  "
$ws.Range("B27").Value = "'EthereumClassic"
$ws.Range("C27").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'17.51"
$ws.Range("E27").Value = "'  +1.90%  "
$ws.Range("B28").Value = "'PancakeSwap"
$ws.Range("C28").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'1.76"
$ws.Range("E28").Value = "'  +13.00%  "
$ws.Range("E29").Value = "'  +3.75%  "
$ws.Range("E30").Value = "'  +0.11%  "
$ws.Range("D31").Value = "'0.0555"
$ws.Range("E31").Value = "'  +1.42%  "
$ws.Range("D32").Value = "'3.97"
$ws.Range("E32").Value = "'  -0.96%  "
$ws.Range("D33").Value = "'3.94"
$ws.Range("E33").Value = "'  -0.87%  "
$ws.Range("E34").Value = "'  +23.31%  "
$ws.Range("E35").Value = "'  +10.85%  "
$ws.Range("E36").Value = "'  -2.69%  "
$ws.Range("D37").Value = "'0.756"
$ws.Range("E37").Value = "'  +7.70%  "
$ws.Range("E38").Value = "'  +10.21%  "
$ws.Range("B39").Value = "'VeChain"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0200"
$ws.Range("E39").Value = "'  +4.43%  "
$ws.Range("B40").Value = "'Aave"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'89.99"
$ws.Range("E40").Value = "'  -1.69%  "
$ws.Range("D41").Value = "'1.346.19"
$ws.Range("E41").Value = "'  +2.10%  "
$ws.Range("D42").Value = "'14.52"
$ws.Range("E42").Value = "'  +1.49%  "
$ws.Range("E43").Value = "'  +2.78%  "
$ws.Range("B44").Value = "'MXToken"
$ws.Range("C44").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.76"
$ws.Range("E44").Value = "'  +4.00%  "
$ws.Range("B45").Value = "'HuobiToken"
$ws.Range("C45").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "'2.37"
$ws.Range("E45").Value = "'  -3.95%  "
$ws.Range("D46").Value = "'11.78"
$ws.Range("E46").Value = "'  +78.24%  "
$ws.Range("D47").Value = "'0.0530"
$ws.Range("E47").Value = "'  +3.69%  "
$ws.Range("D48").Value = "'6.32"
$ws.Range("E48").Value = "'  +1.81%  "
$ws.Range("D49").Value = "'2.028.16"
$ws.Range("E49").Value = "'  +1.50%  "
$ws.Range("E50").Value = "'  +15.44%  "
$ws.Range("D51").Value = "'0.0672"
$ws.Range("E51").Value = "'  +0.44%  "
